# WIP commit for machine migration
# Rebuild the "suite-demo1" sheet (sheet2.xml) as three stacked mini test
# tables: testFoo, testRetailHomepage, testEmailServices - replacing the
# single testEmailServices table that was there before.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("suite-demo1")

# Start clean - remove the old single-table content (old A1:D6 block).
$ws.Cells.Clear()

# ---- helper-ish inline blocks -------------------------------------------
# Each "table" has a bold header row (Function/Target/Value under the test
# name) followed by plain data rows. Column C already carries a left-
# aligned column style and column D a general numeric column style, so
# those are picked up automatically from the `<cols>` definitions; only the
# header row needs explicit bold formatting.

# ---- Table 1: testFoo (rows 2-4) ----------------------------------------
$ws.Range("A2").Value = "testFoo"
$ws.Range("A2").Font.Bold = $true

$ws.Range("B2").Value = "Function"
$ws.Range("B2").Font.Bold = $true

$ws.Range("C2").Value = "Target"
$ws.Range("C2").Font.Bold = $true
$ws.Range("C2").HorizontalAlignment = -4131

$ws.Range("D2").Value = "Value"
$ws.Range("D2").Font.Bold = $true

$ws.Range("B3").Value = "open"
$ws.Range("C3").Value = "/"

$ws.Range("B4").Value = "matt"
$ws.Range("C4").Value = "`$message"

# ---- Table 2: testRetailHomepage (rows 6-12) -----------------------------
$ws.Range("A6").Value = "testRetailHomepage"
$ws.Range("A6").Font.Bold = $true

$ws.Range("B6").Value = "Function"
$ws.Range("B6").Font.Bold = $true

$ws.Range("C6").Value = "Target"
$ws.Range("C6").Font.Bold = $true
$ws.Range("C6").HorizontalAlignment = -4131

$ws.Range("D6").Value = "Value"
$ws.Range("D6").Font.Bold = $true

$ws.Range("B7").Value = "open"
$ws.Range("C7").Value = "/domainname/"

$ws.Range("B8").Value = "click"
$ws.Range("C8").Value = "id=`$searchBox"

$ws.Range("B9").Value = "type"
$ws.Range("C9").Value = "id=`$searchBox"
$ws.Range("D9").Value = "`$domainSearch"

$ws.Range("B10").Value = "click"
$ws.Range("C10").Value = "id=`$searchButton"

$ws.Range("B11").Value = "waitForPageToLoad"
$ws.Range("C11").Value = 10
$ws.Range("C11").NumberFormat = "0"

$ws.Range("B12").Value = "assertTextPresent"
$ws.Range("C12").Value = "`$domainName is available"

# ---- Table 3: testEmailServices (rows 14-16) -----------------------------
$ws.Range("A14").Value = "testEmailServices"
$ws.Range("A14").Font.Bold = $true

$ws.Range("B14").Value = "Function"
$ws.Range("B14").Font.Bold = $true

$ws.Range("C14").Value = "Target"
$ws.Range("C14").Font.Bold = $true
$ws.Range("C14").HorizontalAlignment = -4131

$ws.Range("D14").Value = "Value"
$ws.Range("D14").Font.Bold = $true

$ws.Range("B15").Value = "open"
$ws.Range("C15").Value = "/email-services/"

$ws.Range("B16").Value = "matt"
$ws.Range("C16").Value = "email"

# ---- Trailing styled-but-empty cell (mirrors the old C6 placeholder) ----
$ws.Range("C19").NumberFormat = "0"
